$d = $word.ActiveDocument

# --- Hunk 1: drop the "Meta description: Discover the paying symbols..."
#     paragraph that currently sits right after the H1 title paragraph. ---
$metaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "Meta description*") {
        $metaIndex = $i
        break
    }
}
if ($metaIndex -ge 1) {
    $d.Paragraphs.Item($metaIndex).Range.Delete()
}

# --- Hunk 2: the trailing "Create a feature image for ..." paragraph gets a
#     new bold paragraph inserted right before it ("Play Dragon's Temple for
#     Free - Read Our Review"), and its own (italic) text is swapped for the
#     meta-description copy that used to live near the top. ---
$imgIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "Create a feature image*") {
        $imgIndex = $i
        break
    }
}

$imgPara = $d.Paragraphs.Item($imgIndex)
$insertPoint = $d.Range($imgPara.Range.Start, $imgPara.Range.Start)
$insertPoint.InsertAfter("Play Dragon's Temple for Free - Read Our Review`r")

$newHeadingPara = $d.Paragraphs.Item($imgIndex)
$newHeadingRange = $d.Range($newHeadingPara.Range.Start, $newHeadingPara.Range.End - 1)
$newHeadingRange.Font.Bold = $true

$finalPara = $d.Paragraphs.Item($imgIndex + 1)
$finalTextRange = $d.Range($finalPara.Range.Start, $finalPara.Range.End - 1)
$finalTextRange.Text = "Discover the paying symbols and unique features of Dragon's Temple, including its Chinese theme and music. Play for free and read our review today."
